# Update automatico via Actualizar 12-15-2020 17-13-35
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UF_IVP_DIARIO")

# New rows of data to append (dates are Excel serial numbers, UF and IVP values)
$newData = @(
    @(44197, 29069.39, 30231.51),
    @(44198, 29068.46, 30233.73),
    @(44199, 29067.52, 30235.95),
    @(44200, 29066.58, 30238.17),
    @(44201, 29065.64, 30240.39),
    @(44202, 29064.7,  30242.61),
    @(44203, 29063.759999999998, 30244.83),
    @(44204, 29062.83, 30247.05),
    @(44205, 29061.89, 30249.27)
)

$startRow = 735
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $vals = $newData[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $vals[0]
    $cellA.NumberFormat = "m/d/yyyy"

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $vals[1]
    $cellB.NumberFormat = "0.0000"

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.Value = $vals[2]
    $cellC.NumberFormat = "0.0000"
}

$lastRow = $startRow + $newData.Count - 1

# Update the named range to cover the new data extent
$wb.Names.Item("UF_IVP_DIARIO").RefersTo = "=UF_IVP_DIARIO!`$A`$1:`$C`$" + $lastRow

# Update the selected / active cell to reflect the new bottom of the data
# (frozen panes at row 3 / column 1 are already in place from the source file)
$ws.Activate()
$ws.Cells.Item($lastRow, 2).Select()
